# Add "Oracle session" to the calendar:
# Replace the "06.03: No lecture" entry with "06.03: Oracle session in Aud J"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds the lecture/session description for each week.
# Row 10 currently reads "06.03: No lecture" - update it to the new Oracle session text.
$ws.Range("D10").Value = "06.03: <strong>Oracle session</strong> in Aud J"

# Update the active selection to match the edited cell's new location (D11 after insert-like shift)
$ws.Range("D11").Select()

$wb.Save()
